$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Physical Card Printing" row: cost + reason for the prototype print run
$ws.Range("C8").Value = 210
$ws.Range("D8").Value = "This is just for the intial prototype run of 18 decks of 30 cards for use in showcasing and testing"

# Update convention travel/accommodation estimate
$ws.Range("C13").Value = 500

# Add Total / Rounded labels next to the sum
$ws.Range("G28").Value = "Total"

$ws.Range("F29").Value = 5000
$ws.Range("F29").NumberFormat = "_-""£""* #,##0.00_-;\-""£""* #,##0.00_-;_-""£""* ""-""??_-;_-@_-"
$ws.Range("G29").Value = "Rounded"

# Move active selection like the author left it
$ws.Range("H17").Select() | Out-Null

$wb.Application.CalculateFull() | Out-Null
